# menu item deletion when branch is closed
# Remove every menu row whose "branch" column (D) equals "NTU",
# since that branch is closed. Deleting full rows shifts the rows
# below upward, which matches the target workbook layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$branchColumn = 4   # column D = "branch"
$lastDataRow = 12   # last row of the menu table before the edit

for ($row = $lastDataRow; $row -ge 2; $row--) {
    $branch = $ws.Cells.Item($row, $branchColumn).Value2
    if ($branch -eq "NTU") {
        $ws.Rows.Item($row).Delete()
    }
}
